$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("survey")

$ws.Range("C12").Value = "integer"
$ws.Range("E12").Value = "FA_within_five_meters"
$ws.Range("F12").Value = "Within 5 meters"

$ws.Range("C13").Value = "integer"
$ws.Range("E13").Value = "FA_closest_to_focal"
$ws.Range("F13").Value = "Closest to focal"

$ws.Range("C14").Select()
